$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 237.33333
$ws.Range("I5").Value = 327.75
$ws.Range("K5").Value = 327.75
$ws.Range("M5").Value = -212.75

$ws.Range("H12").Value = 132.55556
$ws.Range("I12").Value = 126
$ws.Range("J12").Value = 135.83333
$ws.Range("K12").Value = 126
$ws.Range("L12").Value = 135.83333
$ws.Range("M12").Value = 44
$ws.Range("N12").Value = -475.83333

$ws.Range("H33").Value = 89.375
$ws.Range("I33").Value = 93.57143000000001
$ws.Range("K33").Value = 93.57143000000001
$ws.Range("M33").Value = 135.42857

$ws.Range("H88").Value = 3999.6
$ws.Range("I88").Value = 3999.5
$ws.Range("J88").Value = 3999.6667
$ws.Range("K88").Value = 3999.5
$ws.Range("L88").Value = 3999.6667
$ws.Range("M88").Value = -3593.5
$ws.Range("N88").Value = -4811.6667

$ws.Range("H91").Value = 3999.6
$ws.Range("I91").Value = 3999.5
$ws.Range("J91").Value = 3999.6667
$ws.Range("K91").Value = 3999.5
$ws.Range("L91").Value = 3999.6667
$ws.Range("M91").Value = -2595.5
$ws.Range("N91").Value = -6807.6667

$ws.Range("H100").Value = 4919.846
$ws.Range("I100").Value = 3067.7144
$ws.Range("K100").Value = 3067.7144
$ws.Range("M100").Value = -2526.7144

$ws.Range("H116").Value = 2295.8572
$ws.Range("I116").Value = 2412
$ws.Range("K116").Value = 2412
$ws.Range("M116").Value = 1030

$ws.Range("H121").Value = 1169.5
$ws.Range("J121").Value = 1169.5
$ws.Range("L121").Value = 3508.5
$ws.Range("N121").Value = -7002.5

$ws.Range("H135").Value = 2375.0908
$ws.Range("I135").Value = 2199.3333
$ws.Range("K135").Value = 19793.9997
$ws.Range("M135").Value = -17258.9997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2147.1428
$ws.Range("I2").Value = 1775
$ws.Range("J2").Value = 3338
$ws.Range("K2").Value = 1775
$ws.Range("L2").Value = 3338
$ws.Range("M2").Value = -1662
$ws.Range("N2").Value = -3564

$ws.Range("H110").Value = 5300
$ws.Range("I110").Value = 4200
$ws.Range("K110").Value = 4200
$ws.Range("M110").Value = -2155

$ws.Range("H116").Value = 2147.1428
$ws.Range("I116").Value = 1775
$ws.Range("J116").Value = 3338
$ws.Range("K116").Value = 1775
$ws.Range("L116").Value = 3338
$ws.Range("M116").Value = 519
$ws.Range("N116").Value = -7926

$ws.Range("H132").Value = 2263.818
$ws.Range("I132").Value = 1244.125
$ws.Range("K132").Value = 3732.375
$ws.Range("M132").Value = -1202.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2147.1428
$ws.Range("I3").Value = 1775
$ws.Range("J3").Value = 3338
$ws.Range("K3").Value = 1775
$ws.Range("L3").Value = 3338
$ws.Range("M3").Value = -1661
$ws.Range("N3").Value = -3566

$ws.Range("H107").Value = 3240.3333
$ws.Range("I107").Value = 1221.6316
$ws.Range("K107").Value = 1221.6316
$ws.Range("M107").Value = 698.3684000000001

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 150000
$ws.Range("J132").Value = 150000
$ws.Range("L132").Value = 150000
$ws.Range("N132").Value = -160120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 67252.86
$ws.Range("J68").Value = 67252.86
$ws.Range("L68").Value = 67252.86
$ws.Range("N68").Value = -68750.86

$ws.Range("H71").Value = 67252.86
$ws.Range("J71").Value = 67252.86
$ws.Range("L71").Value = 201758.58
$ws.Range("N71").Value = -209246.58

$ws.Range("H107").Value = 471.0909
$ws.Range("I107").Value = 453.8
$ws.Range("J107").Value = 644
$ws.Range("K107").Value = 453.8
$ws.Range("L107").Value = 644
$ws.Range("M107").Value = 1466.2
$ws.Range("N107").Value = -4484

$ws.Range("H134").Value = 1497.8649
$ws.Range("I134").Value = 1323.4445
$ws.Range("K134").Value = 3970.3335
$ws.Range("M134").Value = -1435.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2740.6667
$ws.Range("I69").Value = 2222
$ws.Range("K69").Value = 6666
$ws.Range("M69").Value = -5855

$ws.Range("H72").Value = 2740.6667
$ws.Range("I72").Value = 2222
$ws.Range("K72").Value = 19998
$ws.Range("M72").Value = -15942

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2400
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 2333.3333
$ws.Range("K80").Value = 2500
$ws.Range("L80").Value = 2333.3333
$ws.Range("M80").Value = -1502
$ws.Range("N80").Value = -4329.3333

$ws.Range("H83").Value = 2400
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 2333.3333
$ws.Range("K83").Value = 12500
$ws.Range("L83").Value = 11666.6665
$ws.Range("M83").Value = -7508
$ws.Range("N83").Value = -21650.6665

$ws.Range("H107").Value = 462
$ws.Range("I107").Value = 462
$ws.Range("K107").Value = 462
$ws.Range("M107").Value = 1458

$ws.Range("H122").Value = 203399.31
$ws.Range("I122").Value = 315393.5
$ws.Range("J122").Value = 4298.5557
$ws.Range("K122").Value = 946180.5
$ws.Range("L122").Value = 12895.6671
$ws.Range("M122").Value = -943730.5
$ws.Range("N122").Value = -17795.6671

$ws.Range("H128").Value = 34000
$ws.Range("J128").Value = 34000
$ws.Range("L128").Value = 34000
$ws.Range("N128").Value = -43960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2901.8215
$ws.Range("I61").Value = 2110.913
$ws.Range("J61").Value = 6540
$ws.Range("K61").Value = 2110.913
$ws.Range("L61").Value = 6540
$ws.Range("M61").Value = -1908.913
$ws.Range("N61").Value = -6944

$ws.Range("H113").Value = 2901.8215
$ws.Range("I113").Value = 2110.913
$ws.Range("J113").Value = 6540
$ws.Range("K113").Value = 2110.913
$ws.Range("L113").Value = 6540
$ws.Range("M113").Value = 59.08699999999999
$ws.Range("N113").Value = -10880

$ws.Range("H132").Value = 4303.7393
$ws.Range("I132").Value = 3761.2856
$ws.Range("K132").Value = 11283.8568
$ws.Range("M132").Value = -8753.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8331.444
$ws.Range("I62").Value = 5662.6665
$ws.Range("J62").Value = 9665.833000000001
$ws.Range("K62").Value = 5662.6665
$ws.Range("L62").Value = 9665.833000000001
$ws.Range("M62").Value = -5038.6665
$ws.Range("N62").Value = -10913.833

$ws.Range("H65").Value = 8331.444
$ws.Range("I65").Value = 5662.6665
$ws.Range("J65").Value = 9665.833000000001
$ws.Range("K65").Value = 28313.3325
$ws.Range("L65").Value = 48329.165
$ws.Range("M65").Value = -25193.3325
$ws.Range("N65").Value = -54569.165

Write-Host "Applied all profit updates"

